$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the auto-generated chart-tracking defined names
#    (_xlchart.v1.0 / _xlchart.v1.1) that Excel drops once the chart is
#    repositioned / recreated.
# ---------------------------------------------------------------------------
while ($wb.Names.Count -gt 0) {
    $wb.Names.Item(1).Delete()
}

# ---------------------------------------------------------------------------
# 2. New column widths for the "features analysis" block (M, N, Q)
# ---------------------------------------------------------------------------
$ws.Columns.Item(13).ColumnWidth = 17.28515625
$ws.Columns.Item(14).ColumnWidth = 13.42578125
$ws.Columns.Item(17).ColumnWidth = 14

# ---------------------------------------------------------------------------
# 3. Populate the new "Анализ признаков в датасете" table (columns M:Q)
#    Values are entered in the exact order needed so new shared-string
#    entries come out in the same order as the source workbook.
# ---------------------------------------------------------------------------
$ws.Range("N4").Value = "Ср. значение"
$ws.Range("O4").Value = "Макс"
$ws.Range("P4").Value = "Мин"
$ws.Range("Q4").Value = "Кол-во нулей"
$ws.Range("M3").Value = "Анализ признаков в датасете (размер выборки 1897 примеров)"
$ws.Range("M4").Value = "Признак"

$ws.Range("M5").Value = "foodseats"
$ws.Range("M6").Value = "sportsvenue"
$ws.Range("M7").Value = "servicesnum"
$ws.Range("M8").Value = "museums"
$ws.Range("M9").Value = "parks"
$ws.Range("M10").Value = "theatres"
$ws.Range("M11").Value = "library"
$ws.Range("M12").Value = "cultureorg"
$ws.Range("M13").Value = "musartschool"

$ws.Range("N5").Value = 584.54999999999995
$ws.Range("O5").Value = 38482
$ws.Range("P5").Value = 0
$ws.Range("Q5").Value = 904

$ws.Range("N6").Value = 37.04
$ws.Range("O6").Value = 359
$ws.Range("P6").Value = 0
$ws.Range("Q6").Value = 50

$ws.Range("N7").Value = 44.47
$ws.Range("O7").Value = 730
$ws.Range("P7").Value = 0
$ws.Range("Q7").Value = 786

$ws.Range("N8").Value = 0.56000000000000005
$ws.Range("O8").Value = 35
$ws.Range("P8").Value = 0
$ws.Range("Q8").Value = 1239

$ws.Range("N9").Value = 0.19
$ws.Range("O9").Value = 4
$ws.Range("P9").Value = 0
$ws.Range("Q9").Value = 1581

$ws.Range("N10").Value = 0.19
$ws.Range("O10").Value = 34
$ws.Range("P10").Value = 0
$ws.Range("Q10").Value = 1678

$ws.Range("N11").Value = 2.0299999999999998
$ws.Range("O11").Value = 51
$ws.Range("P11").Value = 0
$ws.Range("Q11").Value = 980

$ws.Range("N12").Value = 2.91
$ws.Range("O12").Value = 53
$ws.Range("P12").Value = 0
$ws.Range("Q12").Value = 537

$ws.Range("N13").Value = 0.86
$ws.Range("O13").Value = 8
$ws.Range("P13").Value = 0
$ws.Range("Q13").Value = 995

# ---------------------------------------------------------------------------
# 4. Apply the same styles used by the analogous existing cells
#    (title -> C3 style, header row -> I3 style, data cells -> I4 style)
# ---------------------------------------------------------------------------
$ws.Range("C3").Copy() | Out-Null
$ws.Range("M3").PasteSpecial(-4122) | Out-Null

$ws.Range("I3").Copy() | Out-Null
$ws.Range("M4:Q4").PasteSpecial(-4122) | Out-Null

$ws.Range("I4").Copy() | Out-Null
$ws.Range("M5:Q13").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 5. Move / resize the chart to its new anchor position
# ---------------------------------------------------------------------------
$co = $ws.ChartObjects().Item(1)
$co.Left = 439.232421875
$co.Top = 270.37496062992125
$co.Width = 630.8037896776575
$co.Height = 278.62503937007875

# ---------------------------------------------------------------------------
# 6. Update the sheet view / selection
# ---------------------------------------------------------------------------
$ws.Range("R18").Select() | Out-Null
